$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 9: time (min) 30 -> 50
$ws.Range("B9").Value = 50

# Update shared string for D9: "ex 1.12" -> "ex 1.12-1.13"
$ws.Range("D9").Value = "ex 1.12-1.13"

# Add new row 10
$ws.Range("A10").Value = 211030
$ws.Range("B10").Value = 90
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = "ex 1.14 and 1.16"

# Update selection to match target
$ws.Range("L19").Select()
